$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.967.40"
$ws.Range("E2").Value = "  +5.06%  "
$ws.Range("D3").Value = "'3.348.64"
$ws.Range("E3").Value = "  +4.84%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'557.07"
$ws.Range("E5").Value = "  +3.60%  "
$ws.Range("D6").Value = "'153.19"
$ws.Range("E6").Value = "  +5.82%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.530"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "'7.50"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").Value = "  +4.58%  "
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("D12").Value = "'3.928.06"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("D16").Value = "'63.002.00"
$ws.Range("E16").Value = "  +5.03%  "
$ws.Range("D17").Value = "'3.352.58"
$ws.Range("E17").Value = "  +5.31%  "
$ws.Range("E18").Value = "  +4.07%  "
$ws.Range("D19").Value = "'13.77"
$ws.Range("E19").Value = "  +5.11%  "
$ws.Range("D20").Value = "'8.44"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "'388.44"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "'0.541"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").Value = "'70.71"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("E25").Value = "  +5.69%  "
$ws.Range("D26").Value = "'8.84"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "'0.0₃0973"
$ws.Range("E27").Value = "  +8.12%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +4.40%  "
$ws.Range("E30").Value = "  +4.23%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'5.62"
$ws.Range("E31").Value = "  +4.55%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'23.05"
$ws.Range("E32").Value = "  +3.03%  "
$ws.Range("D33").Value = "'1.30"
$ws.Range("E33").Value = "  +7.13%  "
$ws.Range("D34").Value = "'6.72"
$ws.Range("E34").Value = "  +2.83%  "
$ws.Range("E35").Value = "  +9.89%  "
$ws.Range("D36").Value = "'159.69"
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("E37").Value = "  +12.13%  "
$ws.Range("D38").Value = "'27.40"
$ws.Range("E38").Value = "  +7.18%  "
$ws.Range("D39").Value = "'0.0747"
$ws.Range("E39").Value = "  +4.85%  "
$ws.Range("D40").Value = "'2.842.23"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("E41").Value = "  +8.41%  "
$ws.Range("D42").Value = "'4.32"
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("D43").Value = "'40.75"
$ws.Range("E43").Value = "  +2.52%  "
$ws.Range("D44").Value = "'0.747"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("E45").Value = "  +3.69%  "
$ws.Range("D46").Value = "'3.394.93"
$ws.Range("E46").Value = "  +4.93%  "
$ws.Range("D47").Value = "'22.09"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").Value = "'6.29"
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").Value = "'0.810"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("D51").Value = "'282.47"
$ws.Range("E51").Value = "  +7.66%  "
